$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B10").Value = "Generic User Guiede"
$ws.Range("C10").Value = "Cortex-M4 Devices"
$ws.Range("D10").Value = "DUI0553.pdf"

$ws.Range("B10:C10").HorizontalAlignment = -4108

$ws.Range("B11").Select()
